$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.356.65'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.801.29'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.575'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '34.61'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.298'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0685'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '2.061.34'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.18'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '1.808.43'
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.638'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '34.354.33'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0788'
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.13'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '170.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.43%  '
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.118'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.99%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0524'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.78'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.82'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.392.28'
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.53'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.671'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0188'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '82.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.41'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.943'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('E45').Value = '  +2.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0511'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.25%  '
$ws.Range('D48').Value = '1.962.99'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.50%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0123'
$ws.Range('E51').Value = '  -3.49%  '
